$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row: merge/restructure header cells A1:D1, remove E1:L1 ---
$ws.Range("A1").Value = "Collector,Team,Cycle,Repayment_collections,Repayment_amount,Pending"
$ws.Range("B1").Value = "Amount,Pending"
$ws.Range("C1").Value = "Amount"
$ws.Range("D1").Value = "Recovery,Talk_time,New_collections,Repayment_new_collections,New_collection_amount_rate,New_collection_count_rate"
$ws.Range("E1:L1").ClearContents()

# --- Data updates ---
$ws.Range("D2").Value = 4
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "1,919,739.00"
$ws.Range("E2").Style = "Normal"
$ws.Range("G2").NumberFormat = "@"
$ws.Range("G2").Value = "1.01"
$ws.Range("G2").Style = "Normal"
$ws.Range("H2").Value = 657
$ws.Range("D3").Value = 4
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "1,607,732.00"
$ws.Range("E3").Style = "Normal"
$ws.Range("G3").NumberFormat = "@"
$ws.Range("G3").Value = "1.06"
$ws.Range("G3").Style = "Normal"
$ws.Range("H3").Value = 65
$ws.Range("H4").Value = 127
$ws.Range("D5").Value = 4
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "1,303,645.00"
$ws.Range("E5").Style = "Normal"
$ws.Range("G5").NumberFormat = "@"
$ws.Range("G5").Value = "0.90"
$ws.Range("G5").Style = "Normal"
$ws.Range("H5").Value = 42
$ws.Range("K5").NumberFormat = "@"
$ws.Range("K5").Value = "5.29"
$ws.Range("K5").Style = "Normal"
$ws.Range("H6").Value = 182
$ws.Range("D7").Value = 5
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "1,506,283.00"
$ws.Range("E7").Style = "Normal"
$ws.Range("G7").NumberFormat = "@"
$ws.Range("G7").Value = "0.78"
$ws.Range("G7").Style = "Normal"
$ws.Range("H7").Value = 300
$ws.Range("K7").NumberFormat = "@"
$ws.Range("K7").Value = "1.50"
$ws.Range("K7").Style = "Normal"
$ws.Range("H8").Value = 150
$ws.Range("H9").Value = 517
$ws.Range("H11").Value = 239
$ws.Range("D12").Value = 5
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "650,300.00"
$ws.Range("E12").Style = "Normal"
$ws.Range("G12").NumberFormat = "@"
$ws.Range("G12").Value = "0.40"
$ws.Range("G12").Style = "Normal"
$ws.Range("H12").Value = 126
$ws.Range("H13").Value = 95
$ws.Range("H14").Value = 229
$ws.Range("D15").Value = 3
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "1,067,224.00"
$ws.Range("E15").Style = "Normal"
$ws.Range("G15").NumberFormat = "@"
$ws.Range("G15").Value = "0.63"
$ws.Range("G15").Style = "Normal"
$ws.Range("H15").Value = 384
$ws.Range("D16").Value = 1
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "12,658.00"
$ws.Range("E16").Style = "Normal"
$ws.Range("G16").NumberFormat = "@"
$ws.Range("G16").Value = "0.01"
$ws.Range("G16").Style = "Normal"
$ws.Range("H16").Value = 478
$ws.Range("H17").Value = 143
$ws.Range("H18").Value = 163

# --- Rename sheet ---
$ws.Name = "repayment_20250913_20250913 (1)"
